$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for rows 2-7 (columns K-T)
$updates = @{
    "K2" = 2
    "L2" = 0.6666666666666666
    "M2" = 0.5205496666666667
    "N2" = 1.561649
    "O2" = 0.07725359616672718
    "P2" = 0.0772535961667272
    "Q2" = 0.04937587104888889
    "R2" = 0.44438283944
    "S2" = 0.07725359616672718
    "T2" = 0.0772535961667272
    "O3" = 0.03061701009865156
    "P3" = 0.03061701009865157
    "S3" = 0.03061701009865156
    "T3" = 0.03061701009865157
    "M4" = 0.7422533333333333
    "N4" = 2.22676
    "O4" = 0.1101561348294152
    "P4" = 0.1101561348294152
    "Q4" = 0.07040520284444444
    "R4" = 0.6336468256
    "S4" = 0.1101561348294152
    "T4" = 0.1101561348294152
    "K5" = 1
    "L5" = 0.3333333333333333
    "M5" = 0.03285233333333334
    "N5" = 0.09855700000000001
    "O5" = 0.00487554032782279
    "P5" = 0.00487554032782279
    "Q5" = 0.003116153324444445
    "R5" = 0.02804537992
    "S5" = 0.00487554032782279
    "T5" = 0.00487554032782279
    "M6" = 4.187059333333333
    "N6" = 12.561178
    "O6" = 0.6213919853887639
    "P6" = 0.6213919853887639
    "Q6" = 0.3971565346311111
    "R6" = 3.574408811679999
    "S6" = 0.6213919853887639
    "T6" = 0.6213919853887639
    "M7" = 1.049175333333333
    "N7" = 3.147526
    "O7" = 0.1557057331886193
    "P7" = 0.1557057331886193
    "Q7" = 0.09951777761777778
    "R7" = 0.89565999856
    "S7" = 0.1557057331886193
    "T7" = 0.1557057331886193
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Updated $($updates.Count) cells"
